$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking text (e.g. "0.999").
# Excel auto-converts such text to a Number on assignment unless the cell
# is pre-formatted as Text ("@"), so mark them first to preserve the
# original inline-string (text) semantics from the source data.
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D12", "D14", "D17", "D19", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D37", "D38", "D40", "D45", "D46", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Cell value updates (per diff) ---
$ws.Range("D2").Value = "42.900.10"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.291.80"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "301.09"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "98.97"
$ws.Range("E6").Value = "  +1.74%  "
$ws.Range("D7").Value = "0.502"
$ws.Range("E7").Value = "  -0.64%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.510"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("D10").Value = "36.03"
$ws.Range("E10").Value = "  +7.12%  "
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "18.49"
$ws.Range("E12").Value = "  +9.56%  "
$ws.Range("E13").Value = "  +1.20%  "
$ws.Range("D14").Value = "6.95"
$ws.Range("D15").Value = "2.646.87"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").Value = "2.268.85"
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").Value = "0.798"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").Value = "42.783.02"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "12.50"
$ws.Range("E19").Value = "  +7.82%  "
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").Value = "67.73"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "235.72"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +10.70%  "
$ws.Range("D25").Value = "1.01"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "24.98"
$ws.Range("E27").Value = "  +2.51%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  +5.54%  "
$ws.Range("D29").Value = "34.48"
$ws.Range("E29").Value = "  +1.55%  "
$ws.Range("D30").Value = "166.56"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "9.12"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "17.62"
$ws.Range("E34").Value = "  +4.52%  "
$ws.Range("D35").Value = "4.65"
$ws.Range("E35").Value = "  -2.62%  "
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("D37").Value = "0.0690"
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.81"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +2.05%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.101"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  +0.71%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("D43").Value = "1.980.05"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").Value = "10.10"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").Value = "17.57"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D48").Value = "55.21"
$ws.Range("E48").Value = "  +4.46%  "
$ws.Range("D49").Value = "2.516.30"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "70.63"
$ws.Range("E51").Value = "  +0.96%  "
